# Generate Report for Archive
#
# Two e2e markdown files ("18645ec6-..." and "f2f02ef3-...") moved from
# "Ready for handoff" to "In Translation". As a side effect of the report
# regeneration the rows for "f2f02ef3-..." and "26bc3644-..." swapped
# places in each per-language worksheet (zh-cn / de-de) and on the
# Overview roll-up sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 3 (18645ec6...): status flips to "In Translation" in both language
# columns (date unchanged).
$ov.Range("E3").Value = "In Translation"
$ov.Range("F3").Value = "In Translation"

# Row 4 used to hold 26bc3644...; it now holds f2f02ef3... (status
# "In Translation", generate date updated).
$ov.Range("A4").Value = "f2f02ef3-93be-492d-81b7-931055812b0f.md"
$ov.Range("B4").Value = "e2e\f2f02ef3-93be-492d-81b7-931055812b0f.md"
$ov.Range("E4").Value = "In Translation"
$ov.Range("F4").Value = "In Translation"
$ov.Range("G4").Value = "2016-08-21 20:52:33"

# Row 5 used to hold f2f02ef3...; it now holds 26bc3644... (status stays
# "Ready for handoff", generate date updated).
$ov.Range("A5").Value = "26bc3644-3c72-48da-8374-903556b24682.md"
$ov.Range("B5").Value = "e2e\26bc3644-3c72-48da-8374-903556b24682.md"
$ov.Range("G5").Value = "2016-08-21 20:51:31"

# The hyperlinks anchored on B4/B5 keep pointing at their original
# targets (rId4 -> 26bc3644's commit, rId5 -> f2f02ef3's commit); only
# the display text follows the new row contents.
foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\f2f02ef3-93be-492d-81b7-931055812b0f.md"
    }
    if ($addr -eq '$B$5') {
        $hl.TextToDisplay = "e2e\26bc3644-3c72-48da-8374-903556b24682.md"
    }
}

# ---------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; H4 = "2016-08-21 20:52:28"; H5 = "2016-08-21 20:51:27"; Suffix = "zh-cn.xlf" },
    @{ Name = "de-de"; H4 = "2016-08-21 20:52:33"; H5 = "2016-08-21 20:51:31"; Suffix = "de-de.xlf" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Row 3 (18645ec6...): status flips to "In Translation".
    $ws.Range("C3").Value = "In Translation"

    # Row 4 used to hold 26bc3644...; it now holds f2f02ef3... (status
    # "In Translation", handoff file/date updated).
    $ws.Range("A4").Value = "f2f02ef3-93be-492d-81b7-931055812b0f.md"
    $ws.Range("C4").Value = "In Translation"
    $ws.Range("G4").Value = "f2f02ef3-93be-492d-81b7-931055812b0f.b1d71b11925597f6637a00955e9603006106c95a." + $lang.Suffix
    $ws.Range("H4").Value = $lang.H4

    # Row 5 used to hold f2f02ef3...; it now holds 26bc3644... (status
    # stays "Ready for handoff", handoff file/date updated).
    $ws.Range("A5").Value = "26bc3644-3c72-48da-8374-903556b24682.md"
    $ws.Range("G5").Value = "26bc3644-3c72-48da-8374-903556b24682.0e1e23bccb03f1598ac46727ecf15b2c2af7428e." + $lang.Suffix
    $ws.Range("H5").Value = $lang.H5

    # As on Overview, the A4/A5 hyperlinks keep their original targets;
    # only the display text follows the swapped row contents.
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$4') {
            $hl.TextToDisplay = "f2f02ef3-93be-492d-81b7-931055812b0f.md"
        }
        if ($addr -eq '$A$5') {
            $hl.TextToDisplay = "26bc3644-3c72-48da-8374-903556b24682.md"
        }
    }
}
